# Applies the "Updated symbol list" data refresh described by the diff:
# for each affected row, the Coin/Link (B/C) columns are occasionally
# rotated to a different coin, and the Price/Volume(1h) (D/E) columns are
# refreshed to new quoted readings.
#
# The source workbook stores B/C/D/E as literal text (not numbers/percents),
# e.g. D2 is the text "306.79", not the number 306.79, and E2 is the text
# "-4.64%", not a percentage-formatted -0.0464. Plain `.Value = "307.15"`
# assignment on a General-formatted cell makes Excel auto-convert numeric-
# looking text to a real number (and "%"-suffixed text to a percentage),
# which would change both the stored type and (for the trailing-zero cases)
# the literal digits. Set-Text forces the cell to Text first so the literal
# string is preserved exactly, then restores the "Normal" style so no extra
# style index lingers on the cell afterwards.

function Set-Text($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
Set-Text $ws.Range("D2") '307.15'
Set-Text $ws.Range("E2") '-4.52%'

# Row 3
Set-Text $ws.Range("D3") '39.39'
Set-Text $ws.Range("E3") '-8.00%'

# Row 4
Set-Text $ws.Range("D4") '5.098'
Set-Text $ws.Range("E4") '-2.46%'

# Row 5
Set-Text $ws.Range("D5") '0.07696'
Set-Text $ws.Range("E5") '-6.54%'

# Row 6
Set-Text $ws.Range("E6") '-1.79%'

# Row 7
Set-Text $ws.Range("D7") '1.636'
Set-Text $ws.Range("E7") '-8.64%'

# Row 8
Set-Text $ws.Range("D8") '0.9157'
Set-Text $ws.Range("E8") '-3.64%'

# Row 9
Set-Text $ws.Range("D9") '0.1019'
Set-Text $ws.Range("E9") '-9.18%'

# Row 10
Set-Text $ws.Range("D10") '0.1746'
Set-Text $ws.Range("E10") '-7.63%'

# Row 11
Set-Text $ws.Range("D11") '0.09332'
Set-Text $ws.Range("E11") '-1.34%'

# Row 12
Set-Text $ws.Range("D12") '0.04426'
Set-Text $ws.Range("E12") '-4.42%'

# Row 13
Set-Text $ws.Range("D13") '0.1056'
Set-Text $ws.Range("E13") '-0.28%'

# Row 14
Set-Text $ws.Range("D14") '0.001252'
Set-Text $ws.Range("E14") '-3.05%'

# Row 15
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-Text $ws.Range("D15") '0.005871'
Set-Text $ws.Range("E15") '3.83%'

# Row 16
$ws.Range("B16").Value = 'UpBots'
$ws.Range("C16").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
Set-Text $ws.Range("D16") '0.007491'
Set-Text $ws.Range("E16") '2,415.57%'

# Row 17
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-Text $ws.Range("D17") '3.362'
Set-Text $ws.Range("E17") '0.01%'

# Row 18
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-Text $ws.Range("D18") '2.433'
Set-Text $ws.Range("E18") '-4.27%'

# Row 19
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-Text $ws.Range("D19") '0.3306'
Set-Text $ws.Range("E19") '-1.82%'

# Row 20
$ws.Range("B20").Value = 'MCDex'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-Text $ws.Range("D20") '6.984'
Set-Text $ws.Range("E20") '-6.29%'

# Row 21
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-Text $ws.Range("D21") '0.1347'
Set-Text $ws.Range("E21") '-2.89%'

# Row 22
$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
Set-Text $ws.Range("D22") '0.2810'
Set-Text $ws.Range("E22") '10.29%'

# Row 23
$ws.Range("B23").Value = 'CoinExToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-Text $ws.Range("D23") '0.04134'
Set-Text $ws.Range("E23") '-0.69%'

# Row 24
Set-Text $ws.Range("E24") '-4.01%'

# Row 25
Set-Text $ws.Range("D25") '0.004103'
Set-Text $ws.Range("E25") '-4.24%'

# Row 26
Set-Text $ws.Range("D26") '0.0001299'
Set-Text $ws.Range("E26") '6.48%'

# Row 38
Set-Text $ws.Range("D38") '0.02433'
Set-Text $ws.Range("E38") '-8.81%'

# Row 39
Set-Text $ws.Range("D39") '0.05190'
Set-Text $ws.Range("E39") '-7.60%'

# Row 40
Set-Text $ws.Range("D40") '0.007906'
Set-Text $ws.Range("E40") '-2.72%'

# Row 41
Set-Text $ws.Range("D41") '0.1320'
Set-Text $ws.Range("E41") '-6.12%'

# Row 42
Set-Text $ws.Range("D42") '0.007152'
Set-Text $ws.Range("E42") '10.34%'

# Row 43
Set-Text $ws.Range("D43") '0.001948'
Set-Text $ws.Range("E43") '-9.27%'

# Row 44
Set-Text $ws.Range("D44") '0.008374'
Set-Text $ws.Range("E44") '9.11%'

# Row 45
Set-Text $ws.Range("D45") '0.3060'
Set-Text $ws.Range("E45") '-12.18%'

# Row 46
Set-Text $ws.Range("D46") '0.00006408'
Set-Text $ws.Range("E46") '-5.15%'

# Row 47
Set-Text $ws.Range("D47") '0.00000000749'
Set-Text $ws.Range("E47") '-0.08%'

# Row 48
Set-Text $ws.Range("D48") '0.002996'
Set-Text $ws.Range("E48") '-26.89%'

# Row 49
Set-Text $ws.Range("D49") '0.004425'
Set-Text $ws.Range("E49") '44.05%'

# Row 50
Set-Text $ws.Range("D50") '0.00002098'
Set-Text $ws.Range("E50") '-0.08%'

# Row 51
Set-Text $ws.Range("D51") '0.0001998'
Set-Text $ws.Range("E51") '-0.08%'

